$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '27.122.01'
$c.Style = $origStyle
$c = $ws.Range('E2')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.08%  '
$c.Style = $origStyle

$c = $ws.Range('D3')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.825.83'
$c.Style = $origStyle
$c = $ws.Range('E3')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.29%  '
$c.Style = $origStyle

$c = $ws.Range('D4')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.015'
$c.Style = $origStyle
$c = $ws.Range('E4')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.55%  '
$c.Style = $origStyle

$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '312.75'
$c.Style = $origStyle
$c = $ws.Range('E5')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.06%  '
$c.Style = $origStyle

$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.011'
$c.Style = $origStyle
$c = $ws.Range('E6')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.35%  '
$c.Style = $origStyle

$c = $ws.Range('D7')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.4634'
$c.Style = $origStyle
$c = $ws.Range('E7')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.47%  '
$c.Style = $origStyle

$c = $ws.Range('D8')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.3636'
$c.Style = $origStyle
$c = $ws.Range('E8')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.34%  '
$c.Style = $origStyle

$c = $ws.Range('D9')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.07311'
$c.Style = $origStyle
$c = $ws.Range('E9')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.16%  '
$c.Style = $origStyle

$c = $ws.Range('D10')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.8745'
$c.Style = $origStyle
$c = $ws.Range('E10')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.70%  '
$c.Style = $origStyle

$c = $ws.Range('D11')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '20.12'
$c.Style = $origStyle
$c = $ws.Range('E11')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.45%  '
$c.Style = $origStyle

$c = $ws.Range('D12')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.888.45'
$c.Style = $origStyle
$c = $ws.Range('E12')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +3.56%  '
$c.Style = $origStyle

$c = $ws.Range('D13')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.07645'
$c.Style = $origStyle
$c = $ws.Range('E13')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +4.59%  '
$c.Style = $origStyle

$c = $ws.Range('D14')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '5.355'
$c.Style = $origStyle
$c = $ws.Range('E14')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -2.03%  '
$c.Style = $origStyle

$c = $ws.Range('D15')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '92.72'
$c.Style = $origStyle
$c = $ws.Range('E15')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.04%  '
$c.Style = $origStyle

$c = $ws.Range('D16')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '6.460'
$c.Style = $origStyle
$c = $ws.Range('E16')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.39%  '
$c.Style = $origStyle

$c = $ws.Range('D17')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.011'
$c.Style = $origStyle
$c = $ws.Range('E17')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.21%  '
$c.Style = $origStyle

$c = $ws.Range('D18')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.000008662'
$c.Style = $origStyle
$c = $ws.Range('E18')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.21%  '
$c.Style = $origStyle

$c = $ws.Range('D19')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.012'
$c.Style = $origStyle
$c = $ws.Range('E19')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.46%  '
$c.Style = $origStyle

$c = $ws.Range('D20')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '27.556.20'
$c.Style = $origStyle
$c = $ws.Range('E20')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +1.43%  '
$c.Style = $origStyle

$c = $ws.Range('D21')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '14.48'
$c.Style = $origStyle
$c = $ws.Range('E21')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.99%  '
$c.Style = $origStyle

$c = $ws.Range('D22')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '5.220'
$c.Style = $origStyle
$c = $ws.Range('E22')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.56%  '
$c.Style = $origStyle

$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '10.57'
$c.Style = $origStyle
$c = $ws.Range('E23')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.93%  '
$c.Style = $origStyle

$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '2.087.57'
$c.Style = $origStyle
$c = $ws.Range('E24')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +1.53%  '
$c.Style = $origStyle

$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.882'
$c.Style = $origStyle
$c = $ws.Range('E25')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.85%  '
$c.Style = $origStyle

$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '151.26'
$c.Style = $origStyle
$c = $ws.Range('E26')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.60%  '
$c.Style = $origStyle

$c = $ws.Range('D27')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '18.44'
$c.Style = $origStyle
$c = $ws.Range('E27')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.50%  '
$c.Style = $origStyle

$c = $ws.Range('D28')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '2.076'
$c.Style = $origStyle
$c = $ws.Range('E28')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -3.64%  '
$c.Style = $origStyle

$c = $ws.Range('D29')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '5.122'
$c.Style = $origStyle
$c = $ws.Range('E29')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -2.76%  '
$c.Style = $origStyle

$c = $ws.Range('D30')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '116.30'
$c.Style = $origStyle
$c = $ws.Range('E30')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.93%  '
$c.Style = $origStyle

$c = $ws.Range('D31')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.08900'
$c.Style = $origStyle
$c = $ws.Range('E31')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.24%  '
$c.Style = $origStyle

$c = $ws.Range('B32')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'HuobiToken'
$c.Style = $origStyle
$c = $ws.Range('C32')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c.Style = $origStyle
$c = $ws.Range('D32')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '2.960'
$c.Style = $origStyle
$c = $ws.Range('E32')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.99%  '
$c.Style = $origStyle

$c = $ws.Range('D33')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.156'
$c.Style = $origStyle
$c = $ws.Range('E33')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.24%  '
$c.Style = $origStyle

$c = $ws.Range('B34')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'Filecoin'
$c.Style = $origStyle
$c = $ws.Range('C34')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.Style = $origStyle
$c = $ws.Range('D34')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '4.480'
$c.Style = $origStyle
$c = $ws.Range('E34')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.33%  '
$c.Style = $origStyle

$c = $ws.Range('B35')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'ImmutableX'
$c.Style = $origStyle
$c = $ws.Range('C35')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.Style = $origStyle
$c = $ws.Range('D35')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.7330'
$c.Style = $origStyle
$c = $ws.Range('E35')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -3.40%  '
$c.Style = $origStyle

$c = $ws.Range('E36')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.35%  '
$c.Style = $origStyle

$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '2.529'
$c.Style = $origStyle
$c = $ws.Range('E37')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +4.68%  '
$c.Style = $origStyle

$c = $ws.Range('D38')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.090'
$c.Style = $origStyle
$c = $ws.Range('E38')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.10%  '
$c.Style = $origStyle

$c = $ws.Range('D39')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.05265'
$c.Style = $origStyle
$c = $ws.Range('E39')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.23%  '
$c.Style = $origStyle

$c = $ws.Range('D40')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.01923'
$c.Style = $origStyle
$c = $ws.Range('E40')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.84%  '
$c.Style = $origStyle

$c = $ws.Range('B41')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'MXToken'
$c.Style = $origStyle
$c = $ws.Range('C41')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c.Style = $origStyle
$c = $ws.Range('D41')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '2.937'
$c.Style = $origStyle
$c = $ws.Range('E41')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -2.10%  '
$c.Style = $origStyle

$c = $ws.Range('B42')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'FraxShare'
$c.Style = $origStyle
$c = $ws.Range('C42')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c.Style = $origStyle
$c = $ws.Range('D42')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '7.259'
$c.Style = $origStyle
$c = $ws.Range('E42')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.58%  '
$c.Style = $origStyle

$c = $ws.Range('D43')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.5215'
$c.Style = $origStyle
$c = $ws.Range('E43')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -2.35%  '
$c.Style = $origStyle

$c = $ws.Range('D44')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.1632'
$c.Style = $origStyle
$c = $ws.Range('E44')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.84%  '
$c.Style = $origStyle

$c = $ws.Range('D45')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '8.295'
$c.Style = $origStyle
$c = $ws.Range('E45')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -2.95%  '
$c.Style = $origStyle

$c = $ws.Range('D46')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.4856'
$c.Style = $origStyle
$c = $ws.Range('E46')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.75%  '
$c.Style = $origStyle

$c = $ws.Range('D47')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '10.38'
$c.Style = $origStyle
$c = $ws.Range('E47')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.62%  '
$c.Style = $origStyle

$c = $ws.Range('D48')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.011'
$c.Style = $origStyle
$c = $ws.Range('E48')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.35%  '
$c.Style = $origStyle

$c = $ws.Range('D49')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '103.74'
$c.Style = $origStyle
$c = $ws.Range('E49')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  +0.00%  '
$c.Style = $origStyle

$c = $ws.Range('D50')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.643'
$c.Style = $origStyle
$c = $ws.Range('E50')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -1.57%  '
$c.Style = $origStyle

$c = $ws.Range('D51')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.06286'
$c.Style = $origStyle
$c = $ws.Range('E51')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '  -0.48%  '
$c.Style = $origStyle
